$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E (data rows 2-51) to Text format so numeric-looking
# strings (e.g. "0.594", "44.564.62") are preserved verbatim as text,
# matching the original inline-string cell content.
$deRange = $ws.Range("D2:E51")
$deRange.NumberFormat = "@"

$ws.Range("D2").Value = '44.564.62'
$ws.Range("E2").Value = '  +3.50%  '
$ws.Range("D3").Value = '2.274.94'
$ws.Range("E3").Value = '  +2.08%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '321.56'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").Value = '106.02'
$ws.Range("E6").Value = '  +5.35%  '
$ws.Range("D7").Value = '0.594'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.572'
$ws.Range("D10").Value = '38.72'
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("D11").Value = '0.0844'
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").Value = '7.91'
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.884'
$ws.Range("E14").Value = '  +1.85%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.624.70'
$ws.Range("E15").Value = '  +2.30%  '
$ws.Range("D16").Value = '14.62'
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("D17").Value = '2.281.81'
$ws.Range("E17").Value = '  +2.66%  '
$ws.Range("D18").Value = '44.376.32'
$ws.Range("E18").Value = '  +3.41%  '
$ws.Range("D19").Value = '14.10'
$ws.Range("E19").Value = '  -6.43%  '
$ws.Range("D20").Value = '0.0000101'
$ws.Range("E20").Value = '  +4.05%  '
$ws.Range("D21").Value = '6.54'
$ws.Range("E21").Value = '  +1.26%  '
$ws.Range("D22").Value = '66.54'
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("D23").Value = '3.22'
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("D24").Value = '239.60'
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  +2.80%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '10.21'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").Value = '38.33'
$ws.Range("E29").Value = '  +11.82%  '
$ws.Range("D30").Value = '6.54'
$ws.Range("E30").Value = '  +1.95%  '
$ws.Range("D31").Value = '164.08'
$ws.Range("E31").Value = '  +4.41%  '
$ws.Range("D32").Value = '20.66'
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("D33").Value = '0.0887'
$ws.Range("E33").Value = '  -2.04%  '
$ws.Range("D34").Value = '2.76'
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("E35").Value = '  +3.04%  '
$ws.Range("D36").Value = '0.117'
$ws.Range("E36").Value = '  +11.00%  '
$ws.Range("E37").Value = '  -1.05%  '
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("D39").Value = '3.99'
$ws.Range("E39").Value = '  +1.59%  '
$ws.Range("D40").Value = '4.48'
$ws.Range("E40").Value = '  -0.52%  '
$ws.Range("D41").Value = '15.68'
$ws.Range("E41").Value = '  +23.18%  '
$ws.Range("D42").Value = '0.0329'
$ws.Range("E42").Value = '  +0.49%  '
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("D44").Value = '1.783.56'
$ws.Range("E44").Value = '  -7.40%  '
$ws.Range("D45").Value = '0.209'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = '87.05'
$ws.Range("E46").Value = '  -2.65%  '
$ws.Range("D47").Value = '5.50'
$ws.Range("E47").Value = '  +1.84%  '
$ws.Range("D48").Value = '60.38'
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").Value = '74.98'
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '1.71'
$ws.Range("E50").Value = '  +4.63%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '104.45'
$ws.Range("E51").Value = '  +1.02%  '

# Restore default (Normal) style on D2:E51 so no stray number-format
# style index lingers on the cells (keeps styles.xml equivalent to original).
$deRange.Style = "Normal"
